$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the displayed precision of the existing yearly data row (D4:P4) ---
# The workbook previously used a custom number format (numFmtId 164, "0.0" -
# one decimal place). The edit switches these cells to the built-in "0.00"
# format (two decimal places) before the new 2020 column is appended, so the
# freshly-pasted Q4 cell picks up the same (already-updated) formatting.
$ws.Range("D4:P4").NumberFormat = "0.00"

# --- Add the 2020 column (Q) ---
# Q3: year header, cloned from the existing P3 (2019) header cell so it keeps
# identical font/border/alignment, then overwrite with the new year value.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$ws.Range("Q3").Value = 2020

# Q4: data value for 2020, cloned from the (already reformatted) P4 cell so
# it matches the rest of the row, then overwrite with the new data value.
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 0.067156049127444606

# --- Reset the stored cursor/selection ---
# The saved file no longer pins the view to a stray "B12" selection; move the
# active cell back to A1 (the sheet's natural default position).
$ws.Range("A1").Select()
